$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("D1").Value = "% Complete"
$ws.Range("E1").Value = "Status"

# Row 2 - Requirements gathering
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Complete"

# Row 3 - Stories creation
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Complete"

# Row 4 - stories sign-off
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "In_Progress"

# Row 6 - Execution
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "Not_Started"

# Row 7 - Deliver
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = "Not_Started"

# Row 8 - Close
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "Not_Started"

# Row 5 - Schedule update
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "On_Hold"

# Apply percentage number format to the "% Complete" column
$ws.Range("D2:D8").NumberFormat = "0%"

# Match the saved selection state
$ws.Range("D4").Select()
